# The sheet held raw county crime data with no column headers.
# Insert a new first row and label the four columns so the data is readable.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at the very top; all existing rows (and the gap at the
# old row 63) shift down by one.
$ws.Rows.Item(1).Insert()

$ws.Range("A1").Value = "County"
$ws.Range("B1").Value = "Index"
$ws.Range("C1").Value = "Violent"
$ws.Range("D1").Value = "Property"

# Move the selection like the author left it positioned after editing.
$ws.Range("E1").Select()

# Page was left set to portrait orientation.
$ws.PageSetup.Orientation = 1
